$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price strings that look numeric (e.g. "232.45")
# but must stay plain text (the workbook stores them as inline strings,
# dotted thousands-separators like "43.341.95" are not valid numbers).
# Force text format before writing so Excel does not reinterpret them,
# then restore the cell style so no stray formatting is introduced.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "43.341.95"
$ws.Range("E2").Value = "  -0.94%  "
Set-TextValue "D3" "2.353.17"
$ws.Range("E3").Value = "  +5.45%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue "D5" "232.45"
Set-TextValue "D6" "0.649"
$ws.Range("E6").Value = "  -0.08%  "
Set-TextValue "D7" "67.21"
$ws.Range("E7").Value = "  +6.85%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +1.98%  "
Set-TextValue "D10" "0.0959"
$ws.Range("E10").Value = "  +0.02%  "
Set-TextValue "D11" "56.80"
$ws.Range("E11").Value = "  -0.10%  "
Set-TextValue "D12" "26.26"
$ws.Range("E12").Value = "  -0.81%  "
Set-TextValue "D13" "2.704.17"
$ws.Range("E13").Value = "  +5.46%  "
$ws.Range("E14").Value = "  -0.75%  "
Set-TextValue "D15" "15.68"
$ws.Range("E15").Value = "  +1.93%  "
Set-TextValue "D16" "6.26"
$ws.Range("E16").Value = "  +2.65%  "
Set-TextValue "D17" "0.841"
$ws.Range("E17").Value = "  +2.34%  "
Set-TextValue "D18" "2.354.98"
$ws.Range("E18").Value = "  +5.44%  "
Set-TextValue "D19" "43.275.42"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("E20").Value = "  -0.44%  "
Set-TextValue "D21" "73.94"
$ws.Range("E21").Value = "  +1.88%  "
Set-TextValue "D22" "6.25"
$ws.Range("E22").Value = "  +4.16%  "
Set-TextValue "D23" "249.29"
$ws.Range("E23").Value = "  +0.48%  "
Set-TextValue "D24" "3.98"
$ws.Range("E24").Value = "  +17.65%  "
$ws.Range("E25").Value = "  +0.05%  "
Set-TextValue "D26" "2.46"
$ws.Range("E26").Value = "  +2.25%  "
Set-TextValue "D27" "2.27"
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("E28").Value = "  +0.60%  "
Set-TextValue "D29" "22.39"
$ws.Range("E29").Value = "  +7.86%  "
Set-TextValue "D30" "172.07"
$ws.Range("E30").Value = "  +1.00%  "
Set-TextValue "D31" "1.54"
$ws.Range("E31").Value = "  +12.13%  "
Set-TextValue "D32" "0.126"
$ws.Range("E32").Value = "  -7.90%  "
$ws.Range("E33").Value = "  +0.15%  "
Set-TextValue "D34" "5.02"
$ws.Range("E34").Value = "  +6.11%  "
Set-TextValue "D35" "0.0693"
$ws.Range("E35").Value = "  -0.27%  "
Set-TextValue "D36" "5.05"
$ws.Range("E36").Value = "  +3.56%  "
$ws.Range("E37").Value = "  +9.73%  "
Set-TextValue "D38" "6.51"
$ws.Range("E38").Value = "  +2.11%  "
Set-TextValue "D39" "3.63"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("E40").Value = "  -1.24%  "
Set-TextValue "D41" "8.97"
$ws.Range("E41").Value = "  +8.99%  "
$ws.Range("E42").Value = "  -0.08%  "
Set-TextValue "D43" "18.25"
$ws.Range("E43").Value = "  +7.07%  "
$ws.Range("E44").Value = "  +9.15%  "
Set-TextValue "D45" "1.22"
$ws.Range("E45").Value = "  +3.02%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "98.54"
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D47" "4.46"
$ws.Range("E47").Value = "  +2.47%  "
$ws.Range("E48").Value = "  +0.61%  "
Set-TextValue "D49" "1.445.80"
$ws.Range("E49").Value = "  +1.20%  "
Set-TextValue "D50" "2.576.08"
$ws.Range("E50").Value = "  +5.52%  "
$ws.Range("E51").Value = "  -1.98%  "
